$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "25.999.54"
Set-TextValue "E2" "  +0.36%  "

Set-TextValue "D3" "1.641.35"
Set-TextValue "E3" "  +0.03%  "

Set-TextValue "E4" "  -0.11%  "

Set-TextValue "D5" "214.90"
Set-TextValue "E5" "  -0.13%  "

Set-TextValue "D6" "0.5092"
Set-TextValue "E6" "  +0.88%  "

Set-TextValue "D7" "1.003"
Set-TextValue "E7" "  -0.14%  "

Set-TextValue "D8" "0.2559"

Set-TextValue "E9" "  -0.67%  "

Set-TextValue "D10" "19.53"
Set-TextValue "E10" "  -0.36%  "

Set-TextValue "D11" "0.07778"

Set-TextValue "D12" "4.279"
Set-TextValue "E12" "  +0.06%  "

Set-TextValue "D13" "1.624.33"
Set-TextValue "E13" "  -2.05%  "

Set-TextValue "D14" "0.5420"
Set-TextValue "E14" "  -0.07%  "

Set-TextValue "D15" "64.13"
Set-TextValue "E15" "  -1.00%  "

Set-TextValue "D16" "0.0₅7681"
Set-TextValue "E16" "  -2.27%  "

Set-TextValue "D17" "26.037.14"
Set-TextValue "E17" "  +0.32%  "

Set-TextValue "D18" "1.003"
Set-TextValue "E18" "  -0.09%  "

Set-TextValue "D19" "198.68"
Set-TextValue "E19" "  +0.21%  "

Set-TextValue "D20" "4.417"
Set-TextValue "E20" "  +0.71%  "

Set-TextValue "D21" "9.895"
Set-TextValue "E21" "  -0.63%  "

Set-TextValue "D22" "6.035"
Set-TextValue "E22" "  +0.97%  "

Set-TextValue "D23" "1.006"
Set-TextValue "E23" "  -0.04%  "

Set-TextValue "D24" "1.870"
Set-TextValue "E24" "  +0.08%  "

Set-TextValue "D25" "141.09"
Set-TextValue "E25" "  +0.79%  "

Set-TextValue "D26" "0.1189"
Set-TextValue "E26" "  +4.21%  "

Set-TextValue "D27" "6.808"
Set-TextValue "E27" "  -0.55%  "

Set-TextValue "D28" "15.63"
Set-TextValue "E28" "  -0.55%  "

Set-TextValue "D29" "1.235"
Set-TextValue "E29" "  -0.65%  "

Set-TextValue "D30" "0.04887"
Set-TextValue "E30" "  -0.73%  "

Set-TextValue "D31" "3.248"
Set-TextValue "E31" "  -0.41%  "

Set-TextValue "D32" "3.159"
Set-TextValue "E32" "  -1.10%  "

Set-TextValue "D33" "1.524"
Set-TextValue "E33" "  -0.55%  "

Set-TextValue "D34" "2.367"
Set-TextValue "E34" "  -0.14%  "

Set-TextValue "D35" "0.9019"

Set-TextValue "D36" "2.584"
Set-TextValue "E36" "  -0.81%  "

Set-TextValue "D37" "1.144.05"
Set-TextValue "E37" "  +0.26%  "

Set-TextValue "D38" "0.5441"
Set-TextValue "E38" "  -1.91%  "

Set-TextValue "D40" "1.003"
Set-TextValue "E40" "  -0.26%  "

Set-TextValue "D41" "2.530"
Set-TextValue "E41" "  -1.14%  "

Set-TextValue "E42" "  +7.38%  "

Set-TextValue "D43" "0.8108"
Set-TextValue "E43" "  -1.25%  "

Set-TextValue "D44" "99.29"
Set-TextValue "E44" "  -0.10%  "

Set-TextValue "D45" "5.404"
Set-TextValue "E45" "  -4.98%  "

Set-TextValue "D46" "1.781.68"

Set-TextValue "D47" "0.4532"
Set-TextValue "E47" "  +0.30%  "

Set-TextValue "B48" "Frax"
Set-TextValue "C48" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D48" "1.003"
Set-TextValue "E48" "  -0.27%  "

Set-TextValue "B49" "Aave"
Set-TextValue "C49" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D49" "55.01"
Set-TextValue "E49" "  -0.49%  "

Set-TextValue "D50" "0.05095"
Set-TextValue "E50" "  +0.43%  "

Set-TextValue "E51" "  -0.41%  "
